{"js": "const replacements = [\n  [\"677\u00d72=1354\", \"684\u00d75=3420\"],\n  [\"178\u00d72=356\", \"773\u00d73=2319\"],\n  [\"212\u00d77=1484\", \"665\u00d77=4655\"],\n  [\"460\u00d79=4140\", \"435\u00d74=1740\"],\n  [\"699\u00d76=4194\", \"563\u00d79=5067\"],\n  [\"458\u00d76=2748\", \"930\u00d76=5580\"],\n  [\"707\u00d75=3535\", \"815\u00d79=7335\"],\n  [\"850\u00d74=3400\", \"461\u00d72=922\"],\n  [\"690\u00d75=3450\", \"255\u00d76=1530\"],\n  [\"930\u00d79=8370\", \"893\u00d75=4465\"],\n  [\"635\u00d72=1270\", \"658\u00d79=5922\"],\n  [\"793\u00d73=2379\", \"177\u00d72=354\"],\n  [\"535\u00d72=1070\", \"457\u00d74=1828\"],\n  [\"337\u00d76=2022\", \"635\u00d73=1905\"],\n  [\"764\u00d77=5348\", \"950\u00d77=6650\"],\n  [\"655\u00d77=4585\", \"720\u00d73=2160\"],\n  [\"426\u00d72=852\", \"888\u00d78=7104\"],\n  [\"537\u00d78=4296\", \"381\u00d78=3048\"],\n  [\"319\u00d74=1276\", \"674\u00d75=3370\"],\n  [\"546\u00d76=3276\", \"386\u00d77=2702\"],\n  [\"512\u00d79=4608\", \"741\u00d76=4446\"],\n  [\"836\u00d74=3344\", \"824\u00d77=5768\"],\n  [\"669\u00d75=3345\", \"420\u00d73=1260\"],\n  [\"813\u00d78=6504\", \"450\u00d79=4050\"],\n  [\"313\u00d74=1252\", \"469\u00d78=3752\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old='677\u00d72=1354'; new='684\u00d75=3420'},\n    @{old='178\u00d72=356'; new='773\u00d73=2319'},\n    @{old='212\u00d77=1484'; new='665\u00d77=4655'},\n    @{old='460\u00d79=4140'; new='435\u00d74=1740'},\n    @{old='699\u00d76=4194'; new='563\u00d79=5067'},\n    @{old='458\u00d76=2748'; new='930\u00d76=5580'},\n    @{old='707\u00d75=3535'; new='815\u00d79=7335'},\n    @{old='850\u00d74=3400'; new='461\u00d72=922'},\n    @{old='690\u00d75=3450'; new='255\u00d76=1530'},\n    @{old='930\u00d79=8370'; new='893\u00d75=4465'},\n    @{old='635\u00d72=1270'; new='658\u00d79=5922'},\n    @{old='793\u00d73=2379'; new='177\u00d72=354'},\n    @{old='535\u00d72=1070'; new='457\u00d74=1828'},\n    @{old='337\u00d76=2022'; new='635\u00d73=1905'},\n    @{old='764\u00d77=5348'; new='950\u00d77=6650'},\n    @{old='655\u00d77=4585'; new='720\u00d73=2160'},\n    @{old='426\u00d72=852'; new='888\u00d78=7104'},\n    @{old='537\u00d78=4296'; new='381\u00d78=3048'},\n    @{old='319\u00d74=1276'; new='674\u00d75=3370'},\n    @{old='546\u00d76=3276'; new='386\u00d77=2702'},\n    @{old='512\u00d79=4608'; new='741\u00d76=4446'},\n    @{old='836\u00d74=3344'; new='824\u00d77=5768'},\n    @{old='669\u00d75=3345'; new='420\u00d73=1260'},\n    @{old='813\u00d78=6504'; new='450\u00d79=4050'},\n    @{old='313\u00d74=1252'; new='469\u00d78=3752'}\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.old\n    $find.Replacement.Text = $r.new\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Execute($r.old, $false, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)\n}\n"}
